# Rename the inline Pearson/BTEC logo pictures that live in the
# document's headers and footers:
#   - headers: BTec_Logo-Orange picture   image2.jpg -> image1.jpg
#   - footers: PearsonLogo.png picture    image1.png -> image2.png
#
# Both headers get the same new name, and both footers get the same new
# name, so we just walk every existing Header/Footer of every Section
# and rename its (single) inline picture.

$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {

    for ($i = 1; $i -le 3; $i++) {
        $hdr = $sec.Headers.Item($i)
        if ($hdr.Exists -and $hdr.Range.InlineShapes.Count -gt 0) {
            for ($j = 1; $j -le $hdr.Range.InlineShapes.Count; $j++) {
                # Re-fetch the shape through its own Range so the handle
                # addresses a single, fresh block (avoids a stale-handle
                # error on some stories).
                $shp = $hdr.Range.InlineShapes.Item($j)
                $shp = $shp.Range.InlineShapes.Item(1)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image1.jpg"
                }
            }
        }
    }

    for ($i = 1; $i -le 3; $i++) {
        $ftr = $sec.Footers.Item($i)
        if ($ftr.Exists -and $ftr.Range.InlineShapes.Count -gt 0) {
            for ($j = 1; $j -le $ftr.Range.InlineShapes.Count; $j++) {
                $shp = $ftr.Range.InlineShapes.Item($j)
                $shp = $shp.Range.InlineShapes.Item(1)
                if ($shp.AlternativeText -like "*PearsonLogo.png") {
                    $shp.Name = "image2.png"
                }
            }
        }
    }
}
